$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.875.64"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.631.51"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "1.636.32"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.565"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +17.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "29.885.55"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "0.0₃0700"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.110"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0486"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").Value = "1.429.56"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  -4.84%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0170"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.58%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.824"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "53.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.65%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.75%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "1.773.65"
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "89.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("E51").Value = "  +5.10%  "
